$wb = $excel.ActiveWorkbook

# Overview sheet: update Status for the 61f88a21 row (row 3) for both zh-cn and de-de columns
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: update Status (B3) and Latest Handoff Datetime (D3) for the 61f88a21 row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-01-19 06:53:00"

# de-de sheet: update Status (B3) and Latest Handoff Datetime (D3) for the 61f88a21 row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-01-19 06:53:10"
